$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New VCC/VREF measurement-series rows appended below the existing table.
# Column A holds the series id, column G the description. Within a couple
# of the pairs the description (G) was typed before the id (A); that order
# is preserved here so new values land in the shared-string table the same
# way they originally did.

# VCC measurement rows
$ws.Range("A72").Value = "VCC01"
$ws.Range("G72").Value = "VCC mit CPU, ohne Kondensator"

$ws.Range("G73").Value = "VCC ohne CPU, ohne Kondensator"
$ws.Range("A73").Value = "VCC02"

$ws.Range("A74").Value = "VCC03"
$ws.Range("G74").Value = "VCC mit CPU, mit Kondensator"

$ws.Range("G75").Value = "VCC ohne CPU, mit Kondensator"
$ws.Range("A75").Value = "VCC04"

# row 76 intentionally left blank

# VREF measurement rows
$ws.Range("G77").Value = "VREF von LM317 3,1V"
$ws.Range("A77").Value = "VREF01"

$ws.Range("A78").Value = "VREF02"
$ws.Range("G78").Value = "VREF von AREF Pin µC"

# Update view to reflect where the user ended up after the edit
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K78").Select()
